$wb = $excel.ActiveWorkbook
$wsSummary = $wb.Worksheets.Item("Summary")
$wsConfig = $wb.Worksheets.Item("Config")

# ------------------------------------------------------------------
# 1. Remove the old "Sheet matches (.*) 3$" row (row 19) from Config.
#    This shifts the rest of the Profit-row block up by one (rows
#    20-24 become 19-23) and drops the two now-unused shared strings.
# ------------------------------------------------------------------
$wsConfig.Rows(19).Delete()

# ------------------------------------------------------------------
# 2. Append the new "Loss row" block (mirrors the Profit row block)
#    starting at row 25 (row 24 stays blank, same as the separator
#    above the Profit row block at row 16).
# ------------------------------------------------------------------
$wsConfig.Range("B25").Value = "Name"
$wsConfig.Range("B25").Font.Bold = $true
$wsConfig.Range("C25").Value = "is"
$wsConfig.Range("D25").Value = "Loss row"
$wsConfig.Range("F25").Value = "Do the same for the loss row, but keep it simple this time"

$wsConfig.Range("B26").Value = "Align"
$wsConfig.Range("B26").Font.Bold = $true
$wsConfig.Range("C26").Value = "is"
$wsConfig.Range("D26").Value = $true

$wsConfig.Range("B27").Value = "Table"
$wsConfig.Range("B27").Font.Bold = $true
$wsConfig.Range("C27").Value = "is"
$wsConfig.Range("D27").Value = "PROFIT_RANGE"

$wsConfig.Range("B28").Value = "Source column value"
$wsConfig.Range("B28").Font.Bold = $true
$wsConfig.Range("C28").Value = "is"
$wsConfig.Range("D28").Value = "Loss"

$wsConfig.Range("B29").Value = "Source row column offset"
$wsConfig.Range("B29").Font.Bold = $true
$wsConfig.Range("C29").Value = "is"
$wsConfig.Range("D29").Value = -1

$wsConfig.Range("B30").Value = "Target table"
$wsConfig.Range("B30").Font.Bold = $true
$wsConfig.Range("C30").Value = "is"
$wsConfig.Range("D30").Value = "PROFIT_RANGE"

$wsConfig.Range("B31").Value = "Target row value"
$wsConfig.Range("B31").Font.Bold = $true
$wsConfig.Range("C31").Value = "is"
$wsConfig.Range("D31").Value = "Loss"

# ------------------------------------------------------------------
# 3. Update the view state: selection on Config moves to D31 (the
#    last cell edited), and the Summary tab becomes the active tab
#    (activated last so it ends up "tabSelected").
# ------------------------------------------------------------------
$wsConfig.Activate()
$wsConfig.Range("D31").Select()

$wsSummary.Activate()
$wsSummary.Range("D8").Select()
